# Correct misleading values in mapping schemes; revised area and cost
# assumptions for all occupancies; revised count assumptions for
# non-residential (sheet: Dwellings_buildings).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# Row, B (classification), C (average_area), D (average_unit_cost),
# F (dwellings_per_building). E (classification_proportion) is unchanged.
$rows = @(
    @{ Row = 2;  B = $null;                                    C = 130;  D = 283.5952768729642;  F = $null },
    @{ Row = 3;  B = "Professional and technical services";    C = 130;  D = 325.7442996742671;  F = $null },
    @{ Row = 4;  B = "All other services";                     C = 130;  D = 294.0024429967427;  F = $null },
    @{ Row = 5;  B = $null;                                    C = 260;  D = 283.5952768729642;  F = 2 },
    @{ Row = 6;  B = "Professional and technical services";    C = 260;  D = 325.7442996742671;  F = 2 },
    @{ Row = 7;  B = "All other services";                     C = 260;  D = 294.0024429967427;  F = 2 },
    @{ Row = 8;  B = $null;                                    C = 450;  D = 283.5952768729642;  F = 3 },
    @{ Row = 9;  B = "Professional and technical services";    C = 450;  D = 325.7442996742671;  F = 3 },
    @{ Row = 10; B = "All other services";                     C = 450;  D = 294.0024429967427;  F = 3 },
    @{ Row = 11; B = $null;                                    C = 900;  D = 283.5952768729642;  F = 5 },
    @{ Row = 12; B = "Professional and technical services";    C = 900;  D = 325.7442996742671;  F = 5 },
    @{ Row = 13; B = "All other services";                     C = 900;  D = 294.0024429967427;  F = 5 },
    @{ Row = 14; B = "Professional and technical services";    C = 1200; D = 426.6938110749186;  F = 5 },
    @{ Row = 15; B = "Professional and technical services";    C = 1200; D = 283.5952768729642;  F = 5 },
    @{ Row = 16; B = "All other services";                     C = 1200; D = 294.0024429967427;  F = 5 },
    @{ Row = 17; B = "Professional and technical services";    C = 3200; D = 426.6938110749186;  F = 10 },
    @{ Row = 18; B = "Professional and technical services";    C = 3200; D = 283.5952768729642;  F = 10 },
    @{ Row = 19; B = "All other services";                     C = 3200; D = 294.0024429967427;  F = 10 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    if ($null -ne $r.B) {
        $ws.Cells.Item($rowNum, 2).Value = $r.B
    }
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    if ($null -ne $r.F) {
        $ws.Cells.Item($rowNum, 6).Value = $r.F
    }
}
